$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 through 66 in column A currently share the same bold/border/centered
# "header" style as row 1 (A1/B1). Reset that styling on A2:A66 so only the
# header row keeps it, leaving the data rows with the default "Normal" style.
$ws.Range("A2:A66").Style = "Normal"
